# Append the latest EUR->ARS quote as a new row at the bottom of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 31

# A31/B31 hold plain date/time strings ("2025-09-21", "15:18:29") that must
# stay literal text rather than being auto-converted to date/time serials by
# Excel's type inference, so force the cell to Text format before writing it.
$ws.Range("A" + $newRow).NumberFormat = "@"
$ws.Range("A" + $newRow).Value = "2025-09-21"
$ws.Range("B" + $newRow).Value = "15:18:29"
$ws.Range("C" + $newRow).Value = "1.00 EUR = 1,777.8410"
